# Edit workbook "locais.xlsx":
#  1. Flip the sign of D19 (longitude) from positive to negative.
#  2. Append a new row (23) describing a "LEGO" location (The LEGO Edimburgh).
#  3. Leave the final selection on D19 (matches the last cell the author touched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix the longitude sign on the existing "Poundland" row.
$ws.Range("D19").Value = -3.2056257210258399

# 2. Add the new LEGO row. Write the shared-string cells in the same order
#    they appear in the source sheet (categoria, icone, descricao, nome) so
#    the generated shared-strings table lines up with the authored workbook.
$ws.Range("A23").Value = "LEGO"
$ws.Range("E23").Value = "lego.png"
$ws.Range("G23").Value = "Horário: 10:00–20:00"
$ws.Range("B23").Value = "The LEGO Edimburgh"
$ws.Range("C23").Value = 51.458988661041197
$ws.Range("D23").Value = -2.5842701752980801

# Match the left-aligned "categoria" column styling used by the rows above it.
$ws.Range("A23").HorizontalAlignment = -4131

# 3. Restore the final selection to D19.
$ws.Range("D19").Select()
